# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" on all sheets
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Let the Status columns auto-fit to the new (longer) text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns (E/F) + Latest HO Xliff Generate Date (G)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2017-01-03 05:25:47"

# zh-cn sheet: Status (C) + Latest Handoff Datetime (H)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2017-01-03 05:25:37"

# de-de sheet: Status (C) + Latest Handoff Datetime (H)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2017-01-03 05:25:47"

# Re-fit the Status columns now that the text is longer (target ~17.216 char
# widths; the host quantizes ColumnWidth to whole screen pixels, so feed it
# the input that lands on the closest reachable notch)
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
